# Updating attendance files 25th feb
# Marks newly-absent participants ("A") on their respective attendance
# columns for the FA-II (B) section sheet. Total Absence (column E) is a
# formula (COUNTIF) that recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cells = @("U12", "T22", "T25", "U26", "U32", "T34", "T36", "T37", "U37", "T38", "U38", "T41", "T43", "T50", "T52", "T53", "T55", "U62", "T66", "T67", "T72", "T76", "T80")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $rng.Value = "A"
    # Matches the formatting Excel applies to the other "A" attendance
    # marks already present on the sheet (explicit WrapText toggles the
    # cell's alignment-applied flag without changing the visible format).
    $rng.WrapText = $False
}
